# "line flows display done" — append a new config row (mumbaiDemand) to the
# dummy config sheet, and refresh the selection/scroll position accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append the new key as a single-cell row (column B is left blank, matching
# the other id placeholders that are populated later by the real pipeline).
$ws.Range("A44").Value = "mumbaiDemand"

# Select the full key column (A1:A44) and scroll back to the top of the
# sheet, matching the refreshed view after the new row was added.
$ws.Range("A1:A44").Select()
